$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1619081784284657,0.1945208102290912,0.672860002190442,0.7034158221868738,1.570796280873402,-1.408888145739562),
    @(0.2074942785077096,0.1934383557349265,0.6758762528771421,0.7014819866207826,1.570796286790826,-1.363302046435018),
    @(0.4117502296915038,0.1885882419516801,0.6893910556777709,0.6928171211837666,1.570796313304808,-1.159046098722391),
    @(0.6991189600253006,0.1817645923503434,0.7084051000618862,0.6806264781953137,1.570796350607462,-0.8716773732721977),
    @(0.9033749112090953,0.176914478567097,0.7219199028625151,0.6719616127582977,1.570796377121443,-0.6674214255595712),
    @(0.9489610112883384,0.1758320240729323,0.7249361535492151,0.6700277771922065,1.570796383038867,-0.6218353262550274),
    @(1.900875918751935,0.05335115204042857,1.199938781730575,0.3175060894984869,1.570796342143612,0.3300795436886814),
    @(1.639788217465948,0.07454629022105426,1.120803451180303,0.3754463736539357,1.570796332944535,0.06899185004407737),
    @(1.361539650321506,0.09713455002854179,1.036466685305622,0.437194977539862,1.570796323140818,-0.2092567089567264),
    @(1.083291083177064,0.1197228098360293,0.9521299194309416,0.4989435814257883,1.570796313337101,-0.4875052679575304),
    @(0.8222033818910772,0.140917948016655,0.8729945888806701,0.5568838655812371,1.570796304138024,-0.748592961602134),
    @(0.5937899691995894,0.1594605826197359,0.8037627911880737,0.6075731089698566,1.570796296090182,-0.977006367608534),
    @(0.410269381593821,0.1743588110469463,0.7481379483041316,0.6482997879435473,1.570796289624078,-1.160526949843115),
    @(0.2789178261977563,0.1850219497284968,0.7083254688063935,0.6774491749365862,1.570796284996081,-1.291878501394849),
    @(0.2004217376457329,0.1913942737993129,0.6845334101078362,0.6948689371597545,1.570796282230376,-1.370374587649487),
    @(0.1672303349600339,0.194088756775214,0.6744731406657185,0.7022347352944619,1.570796281060922,-1.403565989363758),
    @(0.1619081784284635,0.1945208102290914,0.6728600021904421,0.7034158221868742,1.570796280873402,-1.40888814573956)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
